$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -0.1340702663479956
$arr[0,1] = -0.2758992130872696
$arr[0,2] = 0.3670754083147943
$arr[0,3] = 0.4914881146443768
$arr[0,4] = 0.09385790860675949
$arr[0,5] = 0.01768919490591373
$arr[0,6] = 0.2940112530288354
$arr[0,7] = 0.1227914307171113
$arr[0,8] = 0.2759388362258526
$arr[0,9] = -0.001489938197266189
$ws.Range("B2:K2").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -0.293453334447234
$arr[0,1] = 0.3495212869548299
$arr[0,2] = 0.4739339932844123
$arr[0,3] = 0.07630378724679503
$arr[0,4] = 0.0001350735459492769
$arr[0,5] = 0.2764571316688709
$arr[0,6] = 0.1052373093571469
$arr[0,7] = 0.2583847148658881
$arr[0,8] = -0.01904405955723064
$arr[0,9] = -0.182031752916177
$ws.Range("B3:K3").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0.4376286059058361
$arr[0,1] = 0.5620413122354185
$arr[0,2] = 0.1644111061978012
$arr[0,3] = 0.08824239249695551
$arr[0,4] = 0.3645644506198771
$arr[0,5] = 0.1933446283081531
$arr[0,6] = 0.3464920338168943
$arr[0,7] = 0.06906325939377558
$arr[0,8] = -0.09392443396517081
$arr[0,9] = -0.2180070093596886
$ws.Range("B4:K4").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0.7406159457232021
$arr[0,1] = 0.3429857396855849
$arr[0,2] = 0.2668170259847391
$arr[0,3] = 0.5431390841076607
$arr[0,4] = 0.3719192617959367
$arr[0,5] = 0.525066667304678
$arr[0,6] = 0.2476378928815592
$arr[0,7] = 0.0846501995226128
$arr[0,8] = -0.03943237587190501
$arr[0,9] = 0.4767206611340558
$ws.Range("B5:K5").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 1.255012967438235
$arr[0,1] = 1.178844253737389
$arr[0,2] = 1.455166311860311
$arr[0,3] = 1.283946489548587
$arr[0,4] = 1.437093895057328
$arr[0,5] = 1.159665120634209
$arr[0,6] = 0.9966774272752628
$arr[0,7] = 0.8725948518807449
$arr[0,8] = 1.388747888886706
$arr[0,9] = 1.178844253737389
$ws.Range("B6:K6").Value = $arr

$arr = New-Object 'object[,]' 1,9
$arr[0,0] = 0.2348700177716323
$arr[0,1] = 0.5111920758945538
$arr[0,2] = 0.3399722535828299
$arr[0,3] = 0.4931196590915711
$arr[0,4] = 0.2156908846684524
$arr[0,5] = 0.05270319130950599
$arr[0,6] = -0.07137938408501182
$arr[0,7] = 0.444773652920949
$arr[0,8] = 0.2348700177716323
$ws.Range("B7:J7").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 0.5151599734076631
$arr[0,1] = 0.343940151095939
$arr[0,2] = 0.4970875566046802
$arr[0,3] = 0.2196587821815615
$arr[0,4] = 0.0566710888226151
$arr[0,5] = -0.06741148657190271
$arr[0,6] = 0.4487415504340581
$arr[0,7] = 0.2388379152847414
$ws.Range("B8:I8").Value = $arr
$ws.Range("J8:J8").ClearContents()

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.4795802412661804
$arr[0,1] = 0.6327276467749217
$arr[0,2] = 0.3552988723518029
$arr[0,3] = 0.1923111789928565
$arr[0,4] = 0.06822860359833866
$arr[0,5] = 0.5843816406042994
$arr[0,6] = 0.3744780054549828
$ws.Range("B9:H9").Value = $arr
$ws.Range("I9:I9").ClearContents()

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.3919214649192569
$arr[0,1] = 0.1144926904961382
$arr[0,2] = -0.04849500286280822
$arr[0,3] = -0.172577578257326
$arr[0,4] = 0.3435754587486348
$arr[0,5] = 0.1336718235993181
$ws.Range("B10:G10").Value = $arr
$ws.Range("H10:H10").ClearContents()

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 0.0691614752440418
$arr[0,1] = -0.09382621811490459
$arr[0,2] = -0.2179087935094224
$arr[0,3] = 0.2982442434965384
$arr[0,4] = 0.08834060834722172
$ws.Range("B11:F11").Value = $arr
$ws.Range("G11:G11").ClearContents()

$arr = New-Object 'object[,]' 1,4
$arr[0,0] = -0.1606876400509585
$arr[0,1] = -0.2847702154454763
$arr[0,2] = 0.2313828215604846
$arr[0,3] = 0.02147918641116785
$ws.Range("B12:E12").Value = $arr
$ws.Range("F12:F12").ClearContents()

$arr = New-Object 'object[,]' 1,3
$arr[0,0] = -0.3143564178021929
$arr[0,1] = 0.201796619203768
$arr[0,2] = -0.00810701594554874
$ws.Range("B13:D13").Value = $arr
$ws.Range("E13:E13").ClearContents()

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.1836459624741271
$arr[0,1] = -0.02625767267518964
$ws.Range("B14:C14").Value = $arr
$ws.Range("D14:D14").ClearContents()

$arr = New-Object 'object[,]' 1,1
$arr[0,0] = -0.04428949692388896
$ws.Range("B15:B15").Value = $arr
$ws.Range("C15:C15").ClearContents()

$ws.Range("B16:B16").ClearContents()
